$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 738, shifting existing rows 738-837 down to 739-838
$ws.Rows.Item(738).Insert()

# Populate the newly inserted row 738 with the new data record
$ws.Cells.Item(738, 1).Value = 10
$ws.Cells.Item(738, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(738, 3).Value = "La Araucanía"
$ws.Cells.Item(738, 4).Value = 45142
$ws.Cells.Item(738, 5).Value = 9
$ws.Cells.Item(738, 6).Value = "Fruta"
$ws.Cells.Item(738, 7).Value = 100101
$ws.Cells.Item(738, 8).Value = "Berries"
$ws.Cells.Item(738, 9).Value = 100101007
$ws.Cells.Item(738, 10).Value = "Kiwi"
$ws.Cells.Item(738, 11).Value = "Hayward"
$ws.Cells.Item(738, 12).Value = "Especial"
$ws.Cells.Item(738, 13).Value = 128
$ws.Cells.Item(738, 14).Value = 14000
$ws.Cells.Item(738, 15).Value = 14000
$ws.Cells.Item(738, 16).Value = 14000
$ws.Cells.Item(738, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(738, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(738, 19).Value = 1400
$ws.Cells.Item(738, 20).Value = 10
